# Add PDP Draft and Updated Hours Worked
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "1/9 - 1/13" hours for each team member (column E, rows 3-6)
$ws.Range("E3").Value = 2
$ws.Range("E4").Value = 2.5
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# Clear the blank/unused row 7 entry in that column
$ws.Range("E7").Value = ""

# Move the active selection to reflect where editing left off
$ws.Range("I19").Select()
